$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 3: A3 gets "RO.ACT.001.CRE" (no special style) ---
$ws.Range("A3").Style = "Normal"
$ws.Range("A3").Value = "RO.ACT.001.CRE"

# --- Row 4: A4 changes from an empty "Text"-styled cell to a normal cell
#            holding "AD.SEC.014.FON.01" (style cleared back to default) ---
$ws.Range("A4").Style = "Normal"
$ws.Range("A4").Value = "AD.SEC.014.FON.01"

# --- Row 6: A6 becomes an empty "Text"-formatted cell ---
$ws.Range("A6").NumberFormat = "@"

# --- Row 7 (new): A7 empty "Text"-formatted cell, D7 = "AD.SEC.001.FON.99" (Text style) ---
$ws.Range("A7").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "AD.SEC.001.FON.99"

# --- Row 8 (new): A8 empty "Text"-formatted cell only ---
$ws.Range("A8").NumberFormat = "@"

# --- Rows 9-13: add an empty "Text"-formatted A cell alongside existing D cells ---
$ws.Range("A9").NumberFormat = "@"
$ws.Range("A10").NumberFormat = "@"
$ws.Range("A11").NumberFormat = "@"
$ws.Range("A12").NumberFormat = "@"
$ws.Range("A13").NumberFormat = "@"

# --- Rows 14-23 (new): empty "Text"-formatted A cells only ---
$ws.Range("A14").NumberFormat = "@"
$ws.Range("A15").NumberFormat = "@"
$ws.Range("A16").NumberFormat = "@"
$ws.Range("A17").NumberFormat = "@"
$ws.Range("A18").NumberFormat = "@"
$ws.Range("A19").NumberFormat = "@"
$ws.Range("A20").NumberFormat = "@"
$ws.Range("A21").NumberFormat = "@"
$ws.Range("A22").NumberFormat = "@"
$ws.Range("A23").NumberFormat = "@"

# --- Selection: active cell B16 within the B15:B16 block ---
$ws.Range("B15:B16").Select()
